$wb = $excel.ActiveWorkbook

# Rename sheets with updated timestamp-based task order IDs
$wb.Worksheets.Item(1).Name = "GNG_TO-16512555641582"
$wb.Worksheets.Item(2).Name = "NB_TO-16512555667555287"
$wb.Worksheets.Item(3).Name = "RS_TO-16512555667622964"
$wb.Worksheets.Item(4).Name = "TOL_TO-16512555668194845"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16512555669073684"

# Sheet 1 (GNG) - update stim filenames
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-1651255564118198.csv"
$ws1.Range("B3").Value = "GNG_stims-16512555641412036.csv"
$ws1.Range("B4").Value = "go_stims-1651255564142201.csv"
$ws1.Range("B5").Value = "GNG_stims-16512555641571991.csv"

# Sheet 2 (NB) - update stim filenames
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "ZB-match_9-16512555642502027.csv"
$ws2.Range("B3").Value = "TB-16512555661571987.csv"
$ws2.Range("B4").Value = "TB-1651255565712199.csv"
$ws2.Range("B5").Value = "ZB-match_4-1651255564202199.csv"
$ws2.Range("B6").Value = "OB-16512555653111975.csv"
$ws2.Range("B7").Value = "ZB-match_0-1651255564450198.csv"
$ws2.Range("B8").Value = "OB-1651255565199199.csv"
$ws2.Range("B9").Value = "TB-16512555667422378.csv"
$ws2.Range("B10").Value = "OB-16512555652262003.csv"

# Sheet 3 (RS) - swap eyes closed / eyes open
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes open"
$ws3.Range("B3").Value = "eyes closed"

# Sheet 4 (TOL) - update stim filenames
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16512555667877374.csv"
$ws4.Range("B3").Value = "ZM_stims-1651255566763886.csv"
$ws4.Range("B4").Value = "MM_stims-16512555668032181.csv"
$ws4.Range("B5").Value = "ZM_stims-1651255566788739.csv"
$ws4.Range("B6").Value = "MM_stims-16512555668184156.csv"
$ws4.Range("B7").Value = "ZM_stims-1651255566804217.csv"

# Sheet 5 (vSAT) - update stim filenames
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-1651255566849704.csv"
$ws5.Range("B3").Value = "SAT_stims-16512555668246477.csv"
$ws5.Range("B4").Value = "vSAT_stims-16512555668918722.csv"
$ws5.Range("B5").Value = "vSAT_stims-1651255566866415.csv"
